# Ranked.xlsx edit: category editing + elo rating bug fix
#  - Games list loses "The Legend of Zelda: Twilight Princess" (row shifts up)
#  - Several rating numbers across Games/TV Shows/Books corrected
#  - Column widths set for the four label columns
#  - Selection / used range shrink by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column widths (A, D, G, J) ----------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 145/6     # -> stored width 25
$ws.Columns.Item(4).ColumnWidth  = 172/6     # -> stored width 29.5
$ws.Columns.Item(7).ColumnWidth  = 146/6     # -> stored width ~25.1640625
$ws.Columns.Item(10).ColumnWidth = 146/6     # -> stored width ~25.1640625

# ---- Games column (D/E): drop "Twilight Princess", shift rows up -------
$ws.Range("D2").Value = "The Legend of Zelda: Skyward Sword"
$ws.Range("E2").Value = 7

$ws.Range("D3").Value = "The Legend of Zelda: Link's Awakening"
$ws.Range("E3").Value = 3

$ws.Range("D4").Value = "Little Nightmares (DNF)"
$ws.Range("E4").Value = 2

$ws.Range("D5").Value = "The Legend of Zelda: Tears of the Kingdom"
$ws.Range("E5").Value = 3

$ws.Range("D6").Value = "Mario Kart"
$ws.Range("E6").Value = 4

$ws.Range("D7").Value = "Fortnite"
$ws.Range("E7").Value = 3

$ws.Range("D8").Value = "Ori and the Blind Forest"
$ws.Range("E8").Value = 7

# Former row 9 is now empty - clear it out entirely
$ws.Range("D9:E9").ClearContents() | Out-Null

# ---- TV Shows ratings (H) ----------------------------------------------
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 4
$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 7
$ws.Range("H6").Value = 4
$ws.Range("H7").Value = 6
$ws.Range("H8").Value = 3

# ---- Books ratings (K) ---------------------------------------------------
$ws.Range("K2").Value = 3
$ws.Range("K3").Value = 7
$ws.Range("K4").Value = 2
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 1
$ws.Range("K8").Value = 4

# ---- Selection / active cell -------------------------------------------
$ws.Range("E8").Select() | Out-Null
